$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header for column B ("value" -> "first_release_value")
$ws.Range("B1").Value = "first_release_value"

# Full replacement data set for column A (date serials) and column B (values),
# rows 2..84 (83 data points replacing the previous 52).
$dates = @(38398, 38487, 38579, 38671, 38763, 38852, 38944, 39036, 39128, 39217, 39309, 39401, 39493, 39583, 39675, 39767, 39859, 39948, 40040, 40132, 40224, 40313, 40405, 40497, 40589, 40678, 40770, 40862, 40954, 41044, 41136, 41228, 41320, 41409, 41501, 41593, 41685, 41774, 41866, 41958, 42050, 42139, 42231, 42323, 42415, 42505, 42597, 42689, 42781, 42870, 42962, 43054, 43146, 43235, 43327, 43419, 43511, 43600, 43692, 43784, 43876, 43966, 44058, 44150, 44242, 44331, 44423, 44515, 44607, 44696, 44788, 44880, 44972, 45061, 45153, 45245, 45337, 45427, 45519, 45611, 45703, 45792, 45884)
$vals = @(-0.3055163919209463, -0.05930270183239372, 0.6231390923824733, -0.7175835844218028, 0.8366283085903774, 0.1767276857804774, 0.3528382780321806, 0.8301537441992792, -1.966285411910945, 0.7805496083026924, 0.5392191980432131, -0.7807892990731773, 0.3147656593484953, -0.6071410908285912, 0.2561593687207875, -0.2643890380460761, 0.4614603479951001, 0.5995551687457663, -0.88436028068827, -0.1757398580474785, -0.8215943001740271, 0.6, 0.4, 0.6349193672116513, 0.4, -0.5972483967898228, 0.8140368792747665, -0.2012463990342326, 0.4, 0.1, 0.3, -0.340329042122363, 0.8, 0.6, 0.1, -0.2963573102740611, 0.6873521384730878, 0.1052962261794335, 0.7267133658511682, 0.7028875639548886, 0.6036537137213145, 0.1326972526782129, 0.5774359918206358, 0.4228511177569345, 0.421075791701611, 0.1946648784293643, 0.3978154615661396, 0.1561519231779869, 0.3026479405721147, 0.9379874529000176, -0.1353320131999567, 0.009046673504869318, 0.4431994051349051, 0.305320875257479, -0.3, 0.2158019844203096, 1.2, 0.0941018033610419, 0.4323758554038761, 0.01870135355044056, -3.2, -10.91706676560194, 7.25, -0.7999999999999972, -2.400000000000006, 3.480874220397794, 6.400000000000006, -1.099999999999994, 0.4209467346675666, 1.540000000000006, -0.1490000000000009, -0.4999966213670604, -0.7399999999999949, -0.2510000000000048, 0.2000000000000028, 0.1200000000000045, 0.09999999999999432, 0, 0.09999999999999432, 0.3200015876295765, 0.09999771493470178, 0.2000034419242951, 0.03999999999999204)

# Extend the date-format styling (copied from the existing templated row A53)
# down through the newly added rows so the new cells pick up the same
# number format / font / border / alignment as the rest of column A.
$ws.Range("A53").Copy()
$ws.Range("A54:A84").PasteSpecial(-4122)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 2).Value = $vals[$i]
}
